# Update progress (%) values in the Android_UserApp_Status_Tracker sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G3").Value = 100
$ws.Range("G4").Value = 50
$ws.Range("G9").Value = 30
$ws.Range("G14").Value = 40

# Leave the active selection on H14, matching the last cell touched in the session
$ws.Range("H14").Select()
